$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells: force Text format before assignment so that
# numeric-looking strings (e.g. "112.90") are kept verbatim as text,
# matching the original inline-string cell type; the format is reset
# back to Normal immediately after so no stray style persists.
$dCells = @("D2","D3","D5","D6","D7","D10","D11","D12","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D26","D28","D29","D31","D32","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '49.773.36'
$ws.Range("D3").Value = '2.658.71'
$ws.Range("D5").Value = '112.90'
$ws.Range("D6").Value = '327.88'
$ws.Range("D7").Value = '0.526'
$ws.Range("D10").Value = '39.84'
$ws.Range("D11").Value = '19.99'
$ws.Range("D12").Value = '0.0819'
$ws.Range("D14").Value = '7.59'
$ws.Range("D15").Value = '3.074.96'
$ws.Range("D16").Value = '2.649.46'
$ws.Range("D17").Value = '0.867'
$ws.Range("D18").Value = '49.770.07'
$ws.Range("D19").Value = '13.52'
$ws.Range("D21").Value = '6.72'
$ws.Range("D22").Value = '0.0₃0952'
$ws.Range("D23").Value = '269.49'
$ws.Range("D24").Value = '69.29'
$ws.Range("D26").Value = '26.28'
$ws.Range("D28").Value = '10.23'
$ws.Range("D29").Value = '2.22'
$ws.Range("D31").Value = '35.00'
$ws.Range("D32").Value = '49.56'
$ws.Range("D38").Value = '2.05'
$ws.Range("D39").Value = '3.16'
$ws.Range("D40").Value = '24.20'
$ws.Range("D41").Value = '128.26'
$ws.Range("D42").Value = '0.0347'
$ws.Range("D43").Value = '2.30'
$ws.Range("D45").Value = '3.36'
$ws.Range("D46").Value = '2.064.51'
$ws.Range("D47").Value = '2.13'
$ws.Range("D48").Value = '2.21'
$ws.Range("D49").Value = '9.00'
$ws.Range("D50").Value = '5.30'
$ws.Range("D51").Value = '59.44'

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Other cells (Coin name / Link / Volume) -- plain text assignment is safe.
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("E10").Value = '  -3.31%  '
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("E13").Value = '  +2.22%  '
$ws.Range("E14").Value = '  +2.63%  '
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("E23").Value = '  -2.45%  '
$ws.Range("E24").Value = '  -4.45%  '
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +2.00%  '
$ws.Range("E29").Value = '  -0.68%  '
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").Value = '  -5.14%  '
$ws.Range("E32").Value = '  -1.32%  '
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -1.29%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("E40").Value = '  +8.90%  '
$ws.Range("E41").Value = '  +2.93%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E42").Value = '  +8.66%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E43").Value = '  +3.75%  '
$ws.Range("E44").Value = '  -0.61%  '
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("E47").Value = '  +6.81%  '
$ws.Range("E48").Value = '  -2.98%  '
$ws.Range("E49").Value = '  -1.44%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("E51").Value = '  -1.90%  '
